$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 5

$ws.Cells.Item(4, 1).Copy($ws.Cells.Item($row, 1))
$ws.Cells.Item($row, 1).Value2 = 42606.88177083333

$ws.Cells.Item($row, 2).Value2 = 87
$ws.Cells.Item($row, 3).Value2 = 0
$ws.Cells.Item($row, 4).Value2 = 0
$ws.Cells.Item($row, 5).Value2 = 0
$ws.Cells.Item($row, 6).Value2 = 0
$ws.Cells.Item($row, 7).Value2 = 0
$ws.Cells.Item($row, 8).Value2 = 0
$ws.Cells.Item($row, 9).Value2 = 0
$ws.Cells.Item($row, 10).Value2 = 0
$ws.Cells.Item($row, 11).Value2 = 0
$ws.Cells.Item($row, 12).Value2 = 0
$ws.Cells.Item($row, 13).Value2 = 0
$ws.Cells.Item($row, 14).Value2 = "Random"
